# Applies the "more console logs added" data refresh to the Report sheet.
# Updates contract addresses and their associated runtime/gas-used counters
# for rows 2-11 (patients, prescribers, pharmacy).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = "0xdE4bBfF9ea4274d5d367713D6cE89E0C72A34FbB"
$ws.Range("D2").Value = 657
$ws.Range("F2").Value = 754
$ws.Range("H2").Value = 1158
$ws.Range("I2").Value = "'044704"

# --- Row 3 ---
$ws.Range("C3").Value = "0xdE126326483911F550447ee90389354451362d44"
$ws.Range("D3").Value = 632
$ws.Range("F3").Value = 771
$ws.Range("H3").Value = 1439
$ws.Range("I3").Value = "'044704"

# --- Row 4 ---
$ws.Range("C4").Value = "0xa674cbe30032Cfc52e88251d9B41cDb1Ac5E74BB"
$ws.Range("D4").Value = 641
$ws.Range("F4").Value = 672
$ws.Range("H4").Value = 1365

# --- Row 5 ---
$ws.Range("C5").Value = "0x213b30Cd94aB7bEC9042a5aa631ad829F93c4F40"
$ws.Range("D5").Value = 754
$ws.Range("F5").Value = 978
$ws.Range("H5").Value = 1204

# --- Row 6 ---
$ws.Range("C6").Value = "0x4798D35C6F5e3f063e320d0D4400bF3B5c161142"
$ws.Range("D6").Value = 599
$ws.Range("F6").Value = 841
$ws.Range("H6").Value = 1472

# --- Row 7 ---
$ws.Range("C7").Value = "0x6e3b29b1bEbE573c5e1B0B7575e12F0C5E2B2E92"
$ws.Range("D7").Value = 626
$ws.Range("F7").Value = 771
$ws.Range("H7").Value = 1120
$ws.Range("I7").Value = "'044704"

# --- Row 8 (prescriber) ---
$ws.Range("C8").Value = "0x8Ba1Ee5d38dBEB24B34468310b0099De1D80A438"
$ws.Range("S8").Value = 3961
$ws.Range("U8").Value = "'0123978212247821224782"
$ws.Range("W8").Value = 3

# --- Row 9 (prescriber) ---
$ws.Range("C9").Value = "0x5F05cD22c404690DB975b939eA303409F03C8393"
$ws.Range("S9").Value = 1255
$ws.Range("U9").Value = "'01239770"

# --- Row 10 (prescriber) ---
$ws.Range("C10").Value = "0xba542c0b1731bdd33A3a51Fa0523e329052F4782"
$ws.Range("S10").Value = 2500
$ws.Range("U10").Value = "'012397701224782"
$ws.Range("W10").Value = 2

# --- Row 11 (pharmacy) ---
$ws.Range("C11").Value = "0x039F1AaCd0Aa6652f5897a32452c0e2B5b471862"
$ws.Range("Y11").Value = 5200
$ws.Range("Z11").Value = 4852
$ws.Range("AA11").Value = 4344

Write-Host "Applied Reports.xlsx data refresh (more console logs added)"
